$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-19 Tuesday", "2024-11-20 Wednesday"),
    @("27×28=756", "98×34=3332"),
    @("67×76=5092", "96×40=3840"),
    @("98×37=3626", "18×56=1008"),
    @("25×91=2275", "13×50=650"),
    @("32×21=672", "49×23=1127"),
    @("69×51=3519", "66×25=1650"),
    @("97×13=1261", "55×15=825"),
    @("21×53=1113", "41×37=1517"),
    @("24×21=504", "53×13=689"),
    @("67×95=6365", "68×51=3468"),
    @("24×93=2232", "97×12=1164"),
    @("90×88=7920", "87×71=6177"),
    @("45×85=3825", "36×91=3276"),
    @("14×75=1050", "77×73=5621"),
    @("39×17=663", "77×97=7469"),
    @("61×17=1037", "32×94=3008"),
    @("13×66=858", "31×69=2139"),
    @("72×37=2664", "42×42=1764"),
    @("30×79=2370", "15×50=750"),
    @("52×98=5096", "54×51=2754"),
    @("23×39=897", "48×78=3744"),
    @("27×66=1782", "87×57=4959"),
    @("78×63=4914", "68×89=6052"),
    @("86×32=2752", "24×73=1752"),
    @("70×55=3850", "69×28=1932")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
